$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Change Runmode column (C) from "Y" to "N" for all existing data rows (2-25)
$ws.Range("C2:C25").Value = "N"

# Change Results for row 12 from "PASS" to "SKIP"
$ws.Range("D12").Value = "SKIP"

# Add a new test case row (row 26), matching the style used on row 8, which has the same
# per-column style combination that the new row needs (bordered cells, B column wraps text)
$ws.Range("A26").Value = "TestCase_B28"
$ws.Range("B26").Value = "Verify that user is able to sort the documents by TIMES CITED field"
$ws.Range("C26").Value = "Y"
$ws.Range("D26").Value = "PASS"

$ws.Range("A8:D8").Copy()
$ws.Range("A26:D26").PasteSpecial(-4122)

# Update selection to match the final state
$ws.Range("B22").Select()
